$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 64 (pushes existing rows 64+ down by 2)
$ws.Rows("64:65").Insert()

# Row 64: new NPC dialogue row (checkpoint / more rods prompt)
$ws.Cells.Item(64, 1).Value2 = "fractionNpc4"
$ws.Cells.Item(64, 2).Value2 = "fractionNpc"
$ws.Cells.Item(64, 3).Value2 = "Professor of Cartography"
$ws.Cells.Item(64, 4).Value2 = "Well done you saved the student and made it to a checkpoint. Come back to me if you run out of rods."
$ws.Cells.Item(64, 7).Value2 = "Okay"
$ws.Cells.Item(64, 9).Value2 = "More Rods."
$ws.Cells.Item(64, 10).Value2 = "dialogue open @e[tag=fractionNpc] @p fractionNpc5"

# Row 65: new NPC dialogue row (confirm reset)
$ws.Cells.Item(65, 1).Value2 = "farctionNpc5"
$ws.Cells.Item(65, 2).Value2 = "fractionNpc"
$ws.Cells.Item(65, 3).Value2 = "Professor of Cartography "
$ws.Cells.Item(65, 4).Value2 = "Getting more rods will reset the game. `nAre you sure you want to reset the game?"
$ws.Cells.Item(65, 7).Value2 = "No thanks"
$ws.Cells.Item(65, 9).Value2 = "Reset the game. "
$ws.Cells.Item(65, 10).Value2 = "scriptevent fraction:npc 1"

# Fix up the style of the "tp" command cell in the (now shifted) groundskeeper row
# so it matches the rest of the table (the old unique font/style gets dropped).
$ws.Cells.Item(64, 8).Copy()
$ws.Cells.Item(67, 8).PasteSpecial(-4122)
$ws.Cells.Item(67, 8).Value2 = "tp @p 30 96 107 facing 30 96 90"

# Both new rows keep the sheet's default row height (15.75), matching the
# target layout rather than auto-fitting to the wrapped two-line prompt text.
$ws.Rows("64:65").RowHeight = 15.75
